$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 16; $r -le 23; $r++) {
    $ws.Cells.Item($r, 1).Value = $r
}

$ws.Range("B16").Select()
